# Rebuild the small summary table: drop the old "Year" / totals-only table
# and replace it with the fuller "by age and gender" breakdown table
# (dropping the now-empty "By age and gender" data column and removing
# duplicate/empty columns as described in the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out all of the old content/strings first so the shared-string table
# gets rebuilt from scratch with only the strings that are actually used.
$ws.Range("A1").ClearContents()
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("A3").ClearContents()

# Header row
$ws.Range("B1").Value = "Total"
$ws.Range("C1").Value = "Of which: male"
$ws.Range("D1").Value = "Of which: female"
$ws.Range("E1").Value = "By age and gender"
$ws.Range("F1").Value = "age 29 or under"
$ws.Range("G1").Value = "Of which: male"
$ws.Range("H1").Value = "Of which: female"
$ws.Range("I1").Value = "30 to 49"
$ws.Range("J1").Value = "Of which: male"
$ws.Range("K1").Value = "Of which: female"
$ws.Range("L1").Value = "50+"
$ws.Range("M1").Value = "Of which: male"
$ws.Range("N1").Value = "Of which: female"

# 2019 data row
$ws.Range("A2").Value = 2019
$ws.Range("B2").Value = 186
$ws.Range("C2").Value = 101
$ws.Range("D2").Value = 85
$ws.Range("F2").Value = 57
$ws.Range("G2").Value = 31
$ws.Range("H2").Value = 26
$ws.Range("I2").Value = 144
$ws.Range("J2").Value = 97
$ws.Range("K2").Value = 47
$ws.Range("L2").Value = 35
$ws.Range("M2").Value = 23
$ws.Range("N2").Value = 12

# 2018 data row
$ws.Range("A3").Value = 2018
$ws.Range("B3").Value = 196
$ws.Range("C3").Value = 122
$ws.Range("D3").Value = 74
$ws.Range("F3").Value = 43
$ws.Range("G3").Value = 24
$ws.Range("H3").Value = 19
$ws.Range("I3").Value = 134
$ws.Range("J3").Value = 89
$ws.Range("K3").Value = 45
$ws.Range("L3").Value = 19
$ws.Range("M3").Value = 9
$ws.Range("N3").Value = 10
